$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.636.88'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.598.21'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.247'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.50'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '1.822.63'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').Value = '1.623.59'
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.523'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.76'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '26.619.72'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '208.36'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.67%  '
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.47'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.12'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.30'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').Value = '1.275.40'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.619'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -7.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.46'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.839'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.06'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +18.04%  '
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.784'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '64.05'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').Value = '1.735.82'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.41'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.44%  '
